$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Country of the 05 - Bhof-Bots team changes from germany to france
$ws.Range("B9").Value = "france"

# New team block: 03 - Robocops (rows 11-13), copy formatting from the
# existing 05 - Bhof-Bots block (rows 7-9) then overwrite with new values.
$ws.Range("A7:B9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A7:B9").Copy()
$ws.Range("A11").PasteSpecial(-4163)

$ws.Range("A11").Value = "Team"
$ws.Range("B11").Value = "03 – Robocops"
$ws.Range("A12").Value = "School"
$ws.Range("B12").Value = "Schule Birklehof Hinterzarten"
$ws.Range("A13").Value = "Country"
$ws.Range("B13").Value = "belgium"

# New team block: 04 - Schokis (rows 15-17)
$ws.Range("A7:B9").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A7:B9").Copy()
$ws.Range("A15").PasteSpecial(-4163)

$ws.Range("A15").Value = "Team"
$ws.Range("B15").Value = "04 – Schokis"
$ws.Range("A16").Value = "School"
$ws.Range("B16").Value = "HSOG"
$ws.Range("A17").Value = "Country"
$ws.Range("B17").Value = "switzerland"

$ws.Range("B19").Select()
